# QA Roadmap update: widen Selector column, wrap/merge selector cells,
# expand the TC-04 selector text, add the new TC-05 "price format" test
# case block, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Widen column E (Selector) so the longer selector strings are legible.
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 33.63

# ---------------------------------------------------------------------
# 2. Turn on word-wrap for the Selector cells of each existing test case
#    block (TC-01..TC-04 header rows) so long selector text wraps nicely.
# ---------------------------------------------------------------------
$ws.Range("E2").WrapText = $true
$ws.Range("E7").WrapText = $true
$ws.Range("E11").WrapText = $true
$ws.Range("E15").WrapText = $true

# ---------------------------------------------------------------------
# 3. Expand the TC-04 ("Validar carga de productos") selector list to
#    include the inventory-page selectors used by its later steps.
# ---------------------------------------------------------------------
$ws.Range("E15").Value = '[data-test="login-button"], [data-test="username"], [data-test="password"], .inventory_list, .inventory_item, .title, .inventory_item_price'

# ---------------------------------------------------------------------
# 4. Merge the (until now unmerged) Selector column cells that span each
#     4-step test-case block, matching the single logical selector value.
# ---------------------------------------------------------------------
$ws.Range("E7:E10").Merge()
$ws.Range("E11:E14").Merge()
$ws.Range("E15:E18").Merge()

# ---------------------------------------------------------------------
# 5. Add the new TC-05 block (rows 19-22) by cloning the TC-03 block
#    (rows 11-14), which has the same row layout/styles, then overwrite
#    the cell values with the new test case content.
# ---------------------------------------------------------------------
$src = $ws.Range("A11:AA14")
$dst = $ws.Range("A19:AA22")
$src.Copy($dst)

$ws.Range("A19").Value = "TC-05"
$ws.Range("B19").Value = "Validar el formato de precios ($)"
$ws.Range("C19").Value = "P0 (Crítica)"
$ws.Range("D19").Value = "Pendiente"

$ws.Range("E19").Value = '[data-test="login-button"], [data-test="username"], [data-test="password"], .inventory_list, .inventory_item, .title, .inventory_item_price'
$ws.Range("E19").Font.Name = "Arial"
$chars = $ws.Range("E19").Characters(119, 21)
$chars.Font.Bold = $true
$chars.Font.Name = "Arial"

$ws.Range("F19").Value = "✅ SÍ (Playwright)"
$ws.Range("G19").Value = "PASSED 🟢"
$ws.Range("H19").Value = (Get-Date -Year 2026 -Month 2 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = " Navegar a la página"
$ws.Range("K19").Value = "La página carga correctamente."
$ws.Range("L19").Value = "Hecho"

$ws.Range("I20").Value = 2
$ws.Range("J20").Value = " Iniciar sesion con credenciales validas"
$ws.Range("K20").Value = "Los campos de texto aceptan la entrada."
$ws.Range("L20").Value = "Hecho"

$ws.Range("I21").Value = 3
$ws.Range("J21").Value = "Validar carga de la página de inventario."
$ws.Range("K21").Value = 'El título de la página es "Products" y el contenedor de ítems es visible.'
$ws.Range("L21").Value = "Hecho"

$ws.Range("I22").Value = 4
$ws.Range("J22").Value = "Extraer lista de precios y verificar el símbolo de moneda."
$ws.Range("K22").Value = 'Todos los precios contienen el símbolo "$" y el formato es numérico válido.'
$ws.Range("L22").Value = "Hecho"

$ws.Range("E19:E22").Merge()

# ---------------------------------------------------------------------
# 6. Freeze the header row, leaving the active selection on B3.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B3").Select()
